$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.187734603881836
$ws.Range("B1").Value = 2.601245641708374
$ws.Range("C1").Value = 9.328455924987793
$ws.Range("D1").Value = 2.085776805877686
$ws.Range("E1").Value = 1.215054988861084
